$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

# "Advanced" options category currently has, in order:
#   row 98: OptionQuickTestTrigger (Dropdown)
#   row 99: OptionQuickTestNow     (Toggle)
#
# Add a new "OptionDebugOverlay" toggle ahead of those two, so the category
# reads: OptionDebugOverlay, OptionQuickTestTrigger, OptionQuickTestNow.
# Inserting the row shifts the two existing rows down to 99/100 while
# carrying their original formatting with them.
$ws.Rows.Item(98).Insert()

$ws.Range("A98").Value = "CategoryAdvanced"
$ws.Range("B98").Value = "OptionDebugOverlay"
$ws.Range("C98").Value = "Toggle"
$ws.Range("D98").Value = "Off"
$ws.Range("F98").Value = "Show on-screen debug overlay with CSM state"
